$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 - Opschorting aanvaard
$ws.Range("F28").Value = "Opschorting aanvaard"

# Row 29 - Begeleidende brief bij opschorting aanvaard
$ws.Range("F29").Value = "Begeleidende brief bij opschorting aanvaard"

# Row 30 - Opschorting niet aanvaard
$ws.Range("F30").Value = "Opschorting niet aanvaard"

# Row 31 - Begeleidende brief bij opschorting niet aanvaard
$ws.Range("F31").Value = "Begeleidende brief bij opschorting niet aanvaard"

# Row 32 - Schrapping aanvaard
$ws.Range("F32").Value = "Schrapping aanvaard"

# Row 33 - Schrapping niet aanvaard
$ws.Range("F33").Value = "Schrapping niet aanvaard"

# Row 34 - Begeleidende brief bij schrapping niet aanvaard
$ws.Range("F34").Value = "Begeleidende brief bij schrapping niet aanvaard"

# Row 35 - Beroep ingewilligd
$ws.Range("F35").Value = "Beroep ingewilligd"

# Row 36 - Begeleidende brief bij beroep ingewilligd
$ws.Range("F36").Value = "Begeleidende brief bij beroep ingewilligd"

# Row 37 - Beroep onontvankelijk verklaard
$ws.Range("F37").Value = "Beroep onontvankelijk verklaard"

# Row 38 - Begeleidende brief bij beroep onontvankelijk verklaard
$ws.Range("F38").Value = "Begeleidende brief bij beroep onontvankelijk verklaard"

# Row 39 - Beroep verworpen
$ws.Range("F39").Value = "Beroep verworpen"

# Row 40 - Begeleidende brief bij beroep verworpen
$ws.Range("F40").Value = "Begeleidende brief bij beroep verworpen"

# Row 41 - semanticRelation / narrower / narrowerTransitive: rename registratie_attest -> registratieattest
$oldSub = "https://data.omgeving.vlaanderen.be/id/concept/leegstand/stuk/registratie_attest"
$newSub = "https://data.omgeving.vlaanderen.be/id/concept/leegstand/stuk/registratieattest"
foreach ($col in @("L41", "M41", "N41")) {
    $current = $ws.Range($col).Value2
    $updated = $current.Replace($oldSub, $newSub)
    $ws.Range($col).Value = $updated
}

# Row 43 - procedure_inkomend_stuk -> inkomend_procedurestuk
$ws.Range("A43").Value = "https://data.omgeving.vlaanderen.be/id/concept/leegstand/stuk/inkomend_procedurestuk"
$ws.Range("E43").Value = "INKOMEND_PROCEDURESTUK"
$ws.Range("F43").Value = "Inkomend procedurestuk"

# Row 44 - semanticRelation / narrower / narrowerTransitive: update references
foreach ($col in @("L44", "M44", "N44")) {
    $ws.Range($col).Value = "https://data.omgeving.vlaanderen.be/id/concept/leegstand/stuk/inkomend_procedurestuk|https://data.omgeving.vlaanderen.be/id/concept/leegstand/stuk/retour_afzender_procedurestuk"
}

# Row 45 - registratie_attest -> registratieattest
$ws.Range("A45").Value = "https://data.omgeving.vlaanderen.be/id/concept/leegstand/stuk/registratieattest"
$ws.Range("E45").Value = "REGATT"
$ws.Range("F45").Value = "Registratieattest"

# Row 46 - registratie_attest_begeleidende_brief -> registratieattest_begeleidende_brief
$ws.Range("A46").Value = "https://data.omgeving.vlaanderen.be/id/concept/leegstand/stuk/registratieattest_begeleidende_brief"
$ws.Range("C46").Value = "null"
$ws.Range("E46").Value = "REGATT_BEG_BRIEF"
$ws.Range("F46").Value = "Begeleidende brief bij een registratieattest"

# Row 47 - retour_afzender -> retour_afzender_procedurestuk
$ws.Range("A47").Value = "https://data.omgeving.vlaanderen.be/id/concept/leegstand/stuk/retour_afzender_procedurestuk"
$ws.Range("E47").Value = "RETOUR_AFZENDER_PROCEDURESTUK"

# Row 48 - Begeleidende brief bij schrapping aanvaard
$ws.Range("F48").Value = "Begeleidende brief bij schrapping aanvaard"
